# Applies the "processamento de varios .lis" update:
#  - Sheet "Dados": the voltage-interval table is narrowed from 20 rows (bins 23-42)
#    down to 10 rows (bins 31-40), with new frequency / cumulative / percentual values.
#  - Sheet "Estatisticas": mean/variance/std-dev (grouped & ungrouped) plus the computed
#    statistics block (rows 7 and 10-25) are refreshed to match the new dataset.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Dados"
# ---------------------------------------------------------------------------
$dados = $wb.Worksheets.Item("Dados")

# Remove the old trailing rows (12-21) so the table shrinks from A1:F21 to A1:F11
$dados.Range("A12:F21").EntireRow.Delete()

# New interval/frequency table values (rows 2-11)
$dadosValues = @(
    @(2, 31, 1.55, 664424.093, 0,   0,   100),
    @(3, 32, 1.6,  685857.128, 271, 271, 9.666667),
    @(4, 33, 1.65, 707290.1629999999, 19, 290, 3.333333),
    @(5, 34, 1.7,  728723.199, 2,   292, 2.666667),
    @(6, 35, 1.75, 750156.2340000001, 2, 294, 2),
    @(7, 36, 1.8,  771589.269, 0,   294, 2),
    @(8, 37, 1.85, 793022.304, 1,   295, 1.666667),
    @(9, 38, 1.9,  814455.34,  4,   299, 0.333333),
    @(10, 39, 1.95, 835888.375, 0,  299, 0.333333),
    @(11, 40, 2,   857321.41,  1,   300, 0)
)

foreach ($row in $dadosValues) {
    $r = $row[0]
    $dados.Cells.Item($r, 1).Value = $row[1]
    $dados.Cells.Item($r, 2).Value = $row[2]
    $dados.Cells.Item($r, 3).Value = $row[3]
    $dados.Cells.Item($r, 4).Value = $row[4]
    $dados.Cells.Item($r, 5).Value = $row[5]
    $dados.Cells.Item($r, 6).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "Estatisticas"
# ---------------------------------------------------------------------------
$stat = $wb.Worksheets.Item("Estatisticas")

# Mean / Variance / Standard deviation (grouped vs ungrouped), from the .lis files
$stat.Range("B2").Value = 1.586
$stat.Range("C2").Value = 1.59964919
$stat.Range("B3").Value = 0.00220301003
$stat.Range("C3").Value = 0.00193312182
$stat.Range("B4").Value = 0.0469362337
$stat.Range("C4").Value = 0.0439672813

# Computed statistics block (row 7) derived from the refreshed "Dados" table
$stat.Range("A7").Value = 1.611
$stat.Range("B7").Value = 0.002195666666666664
$stat.Range("C7").Value = 0.04685794134046719
$stat.Range("D7").Value = 1.6
$stat.Range("E7").Value = 1.6
$stat.Range("F7").Value = 300
$stat.Range("G7").Value = 0.02908624540066244
$stat.Range("H7").Value = 5.722071434572147
$stat.Range("I7").Value = 35.01731264748563
$stat.Range("J7").Value = 0.2613037285559408

# Long-form metric/value list (rows 10-25) mirroring the blocks above
$stat.Range("B10").Value = 1.586
$stat.Range("B11").Value = 1.59964919
$stat.Range("B12").Value = 0.00220301003
$stat.Range("B13").Value = 0.00193312182
$stat.Range("B14").Value = 0.0469362337
$stat.Range("B15").Value = 0.0439672813
$stat.Range("B16").Value = 1.611
$stat.Range("B17").Value = 0.002195666666666664
$stat.Range("B18").Value = 0.04685794134046719
$stat.Range("B19").Value = 1.6
$stat.Range("B20").Value = 1.6
$stat.Range("B22").Value = 0.02908624540066244
$stat.Range("B23").Value = 5.722071434572147
$stat.Range("B24").Value = 35.01731264748563
$stat.Range("B25").Value = 0.2613037285559408
